$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 166, shifting existing rows 166-178 down to 167-179.
$xlShiftDown = -4121
$ws.Rows("166:166").Insert($xlShiftDown)

# Populate the newly inserted row 166 with the new Chirimoya price record.
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = "Vega Modelo de Temuco"
$ws.Range("C166").Value = "La Araucanía"
$ws.Range("D166").Value = 45127
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100107
$ws.Range("H166").Value = "Otros"
$ws.Range("I166").Value = 100107002
$ws.Range("J166").Value = "Chirimoya"
$ws.Range("K166").Value = "Cultivar IV Región"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 55
$ws.Range("N166").Value = 5000
$ws.Range("O166").Value = 5000
$ws.Range("P166").Value = 5000
$ws.Range("Q166").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R166").Value = "Provincia del Elquí"
$ws.Range("S166").Value = 5000
$ws.Range("T166").Value = 1
